# "Added request status tracker feature":
# Insert a new tracker row just above the first dated time-sheet entry
# (old row 5), pushing the existing entries (old rows 5-13) down to 6-14
# and extending the used range to A1:D14.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("5:5").Insert()

# Materialize empty A5:C5 cells (no date/clock-in/clock-out yet for this
# tracker row) and seed D5 with 0 hours.
$ws.Cells.Item(5, 1).Font.Name = "Calibri"
$ws.Cells.Item(5, 2).Font.Name = "Calibri"
$ws.Cells.Item(5, 3).Font.Name = "Calibri"
$ws.Cells.Item(5, 4).Value = 0
